$d = $word.ActiveDocument

# The existing "_GoBack" bookmark sits at the end of the last paragraph's
# text. In the target revision it is relocated into the middle of the new
# paragraph that gets appended below, so remove it from its old spot first.
$oldBookmark = $d.Bookmarks.Item("_GoBack")
$oldBookmark.Delete()

# Append, at the very end of the document body, a brand-new empty paragraph
# followed by a new paragraph holding the 26/02/2014 log entry (with the
# "_GoBack" bookmark re-inserted mid-sentence, exactly where it used to be
# relative to the trailing text).
$end = $d.Content.End
$insertionRange = $d.Range($end, $end)

$newXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:t xml:space="preserve">26/02/2014 2hrs. Ya puede procesar multiplicaciones. Falta restas y </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>dem&#225;s ,al</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> menos en las constantes globales. </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>GlobalConstants</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> with names that includes previous </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>GlobalConstants</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> as "TAM","TAM2" shouldn't be</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>used since parser recognizes TAM, replaces for its value and then the 2.</w:t></w:r></w:p>
'@

[void]$insertionRange.InsertXML($newXml)
